# Insert a new client "CEDEÑO MACIAS FRANCISCO ARMANDO" as row 14 in both
# the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, pushing the
# alphabetically-later rows (CERAMIKASA S.A.S. onward) down by one row, and
# update the trailing "X de N" summary row on "VENTAS POR GRUPO" to reflect
# the new total client count (53 -> 54).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (17 value columns, C..R)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a blank row before the current row 14 (CERAMIKASA S.A.S.), shifting
# everything below (including the summary row) down by one.
$ws1.Rows.Item(14).Insert()

$ws1.Cells.Item(14, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws1.Cells.Item(14, 2).Value = "CEDEÑO MACIAS FRANCISCO ARMANDO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(14, $col).Value = 0
}

# The summary row (previously row 55, now row 56) carries literal text like
# "0 de 53" / "1 de 53" that must be bumped to reflect 54 total clients.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(56, $col)
    $txt = $cell.Text
    $newTxt = $txt.Replace("de 53", "de 54")
    $cell.Value = $newTxt
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (5 value columns, C..G)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(14).Insert()

$ws2.Cells.Item(14, 1).Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws2.Cells.Item(14, 2).Value = "CEDEÑO MACIAS FRANCISCO ARMANDO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(14, $col).Value = 0
}
